$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.367.73"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "3.372.59"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.79"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.26"
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("D8").Value = "3.363.47"
$ws.Range("E8").Value = "  -1.36%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.631"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.15"
$ws.Range("E12").Value = "  -3.52%  "

$ws.Range("E13").Value = "  -2.53%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "3.922.35"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.18"
$ws.Range("E16").Value = "  -1.07%  "

$ws.Range("D17").Value = "3.389.29"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").Value = "65.408.58"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.79"
$ws.Range("E20").Value = "  -1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "479.62"
$ws.Range("E22").Value = "  +2.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.12"
$ws.Range("E24").Value = "  +4.13%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.29"
$ws.Range("E25").Value = "  +4.73%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.09"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.55"
$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("E29").Value = "  -2.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.09"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.52"
$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "63.68"
$ws.Range("E32").Value = "  +3.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.39"
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "571.96"
$ws.Range("E34").Value = "  -2.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  -1.55%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  +3.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.54"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.372"
$ws.Range("E40").Value = "  -0.79%  "

$ws.Range("D41").Value = "0.0₃0735"
$ws.Range("E41").Value = "  -3.35%  "

$ws.Range("D42").Value = "3.100.04"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.79"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  -1.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -3.52%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.58"
$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.40"
$ws.Range("E51").Value = "  +0.36%  "
